# Add two new columns (I: I0, J: IF) to the sheet, mirroring the style of
# the existing header row and filling in the corresponding numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - copy style from existing header cell (H1) so the new
# headers match the bold/bordered/centered formatting used by the rest
# of row 1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-23
$data = @(
    @(5, 5),
    @(4, 6),
    @(6, 7),
    @(9, 9),
    @(6, 7),
    @(5, 6),
    @(6, 6),
    @(7, 7),
    @(7, 8),
    @(6, 7),
    @(8, 9),
    @(10, 10),
    @(6, 6),
    @(8, 8),
    @(5, 7),
    @(8, 9),
    @(5, 6),
    @(6, 6),
    @(11, 11),
    @(6, 6),
    @(7, 7),
    @(7, 7)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
